$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text so that values like
# "325.66" or "0.4610" are not auto-converted to numbers (which would
# lose trailing zeros / change formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.088.58'
$ws.Range("E2").Value = '  -1.79%  '

$ws.Range("D3").Value = '1.834.85'
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '325.66'
$ws.Range("E5").Value = '  -3.24%  '

$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").Value = '0.4610'
$ws.Range("E7").Value = '  -1.15%  '

$ws.Range("D8").Value = '0.3861'
$ws.Range("E8").Value = '  -1.11%  '

$ws.Range("D9").Value = '0.07852'
$ws.Range("E9").Value = '  -0.70%  '

$ws.Range("D10").Value = '0.9620'
$ws.Range("E10").Value = '  -1.95%  '

$ws.Range("D11").Value = '21.97'
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").Value = '1.861.02'
$ws.Range("E12").Value = '  +1.20%  '

$ws.Range("D13").Value = '5.681'
$ws.Range("E13").Value = '  -2.68%  '

$ws.Range("D14").Value = '6.893'
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("D15").Value = '0.06845'
$ws.Range("E15").Value = '  -0.92%  '

$ws.Range("E16").Value = '  +0.68%  '

$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").Value = '0.000009935'
$ws.Range("E18").Value = '  -0.97%  '

$ws.Range("D19").Value = '16.69'
$ws.Range("E19").Value = '  -2.41%  '

$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.14%  '

$ws.Range("D21").Value = '28.103.82'
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").Value = '5.293'
$ws.Range("E22").Value = '  -1.90%  '

$ws.Range("E23").Value = '  -2.47%  '

$ws.Range("D24").Value = '2.087'
$ws.Range("E24").Value = '  -4.08%  '

$ws.Range("D25").Value = '2.068.29'
$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").Value = '154.57'
$ws.Range("E26").Value = '  +0.72%  '

$ws.Range("D27").Value = '19.12'
$ws.Range("E27").Value = '  -1.72%  '

$ws.Range("D28").Value = '5.744'
$ws.Range("E28").Value = '  -5.52%  '

$ws.Range("D29").Value = '1.974'
$ws.Range("E29").Value = '  -2.79%  '

$ws.Range("D30").Value = '119.25'
$ws.Range("E30").Value = '  +1.37%  '

$ws.Range("D31").Value = '0.9434'
$ws.Range("E31").Value = '  -2.98%  '

$ws.Range("D32").Value = '0.09253'
$ws.Range("E32").Value = '  -1.03%  '

$ws.Range("D33").Value = '5.271'
$ws.Range("E33").Value = '  -1.62%  '

$ws.Range("E34").Value = '  -1.76%  '

$ws.Range("D35").Value = '3.329'
$ws.Range("E35").Value = '  -4.41%  '

$ws.Range("D36").Value = '0.05839'
$ws.Range("E36").Value = '  -5.17%  '

$ws.Range("D37").Value = '0.02114'
$ws.Range("E37").Value = '  -3.86%  '

$ws.Range("D38").Value = '1.136'
$ws.Range("E38").Value = '  -2.50%  '

$ws.Range("D39").Value = '7.711'
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("D40").Value = '0.5598'
$ws.Range("E40").Value = '  -2.06%  '

$ws.Range("D41").Value = '9.906'
$ws.Range("E41").Value = '  -2.43%  '

$ws.Range("D42").Value = '0.1758'
$ws.Range("E42").Value = '  -2.02%  '

$ws.Range("D43").Value = '0.07326'
$ws.Range("E43").Value = '  +3.15%  '

$ws.Range("D44").Value = '11.69'
$ws.Range("E44").Value = '  -0.37%  '

$ws.Range("D45").Value = '0.5266'
$ws.Range("E45").Value = '  -2.23%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.130'
$ws.Range("E46").Value = '  -9.65%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '2.121'
$ws.Range("E47").Value = '  -11.62%  '

$ws.Range("D48").Value = '1.835'
$ws.Range("E48").Value = '  -3.75%  '

$ws.Range("D49").Value = '113.35'
$ws.Range("E49").Value = '  -0.16%  '

$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("D51").Value = '1.021'
$ws.Range("E51").Value = '  +0.04%  '

# Reset the style on column D back to Normal (default), so we don't
# leave a lingering explicit "text" cell style/format behind.
$ws.Range("D2:D51").Style = "Normal"
